# Daily attendance processing - 2026-01-02 11:53:56
# Swap the order of "System" and the recorded-by email address in the
# "Recorded By" column (column G) wherever both are present together.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # Column G = "Recorded By"
    if ($cell.Text -eq $oldValue) {
        $cell.Value = $newValue
    }
}
